# This script renames several shared-string header labels across the
# shortcircuit result-comparison workbook:
#   *_deg        -> *_degree
#   pf_q_*_mw    -> pf_q_*_mvar   (for the phase a/b/c "from"/"to" columns)
#   pf_vm_*_bus_pu -> pf_vm_*_pu  (drop the superfluous "_bus")
#   pf_va_*_bus_deg -> pf_va_*_degree (drop "_bus", expand "_deg")
#
# The three-phase ("LLL_*") sheets only have the simple A1:Q1 header
# (pf_ikss_from_deg, pf_ikss_to_deg, pf_va_from_deg, pf_va_to_deg),
# while every other fault-type sheet (LL_*, LLG_*, LG_*) shares the
# wider A1:AQ1 header that also needs the per-phase columns fixed.

$wb = $excel.ActiveWorkbook

$simpleSheets = @(
    'LLL_max_6', 'LLL_max_10', 'LLL_max_fault_6', 'LLL_max_fault_10',
    'LLL_min_6', 'LLL_min_10', 'LLL_min_fault_6', 'LLL_min_fault_10'
)

$wideSheets = @(
    'LL_max_6', 'LL_max_10', 'LL_max_fault_6', 'LL_max_fault_10',
    'LL_min_6', 'LL_min_10', 'LL_min_fault_6', 'LL_min_fault_10',
    'LLG_max_6', 'LLG_max_10', 'LLG_max_fault_6', 'LLG_max_fault_10',
    'LLG_min_6', 'LLG_min_10', 'LLG_min_fault_6', 'LLG_min_fault_10',
    'LG_max_6', 'LG_max_10', 'LG_max_fault_6', 'LG_max_fault_10',
    'LG_min_6', 'LG_min_10', 'LG_min_fault_6', 'LG_min_fault_10'
)

foreach ($name in $simpleSheets) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("L1").Value = "pf_ikss_from_degree"
    $ws.Range("M1").Value = "pf_ikss_to_degree"
    $ws.Range("P1").Value = "pf_va_from_degree"
    $ws.Range("Q1").Value = "pf_va_to_degree"
}

foreach ($name in $wideSheets) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("T1").Value  = "pf_q_a_from_mvar"
    $ws.Range("U1").Value  = "pf_q_b_from_mvar"
    $ws.Range("V1").Value  = "pf_q_c_from_mvar"
    $ws.Range("W1").Value  = "pf_q_a_to_mvar"
    $ws.Range("X1").Value  = "pf_q_b_to_mvar"
    $ws.Range("Y1").Value  = "pf_q_c_to_mvar"

    $ws.Range("Z1").Value  = "pf_ikss_a_from_degree"
    $ws.Range("AA1").Value = "pf_ikss_b_from_degree"
    $ws.Range("AB1").Value = "pf_ikss_c_from_degree"
    $ws.Range("AC1").Value = "pf_ikss_a_to_degree"
    $ws.Range("AD1").Value = "pf_ikss_b_to_degree"
    $ws.Range("AE1").Value = "pf_ikss_c_to_degree"

    $ws.Range("AG1").Value = "pf_vm_b_from_pu"
    $ws.Range("AH1").Value = "pf_vm_c_from_pu"
    $ws.Range("AI1").Value = "pf_vm_a_to_pu"
    $ws.Range("AJ1").Value = "pf_vm_b_to_pu"
    $ws.Range("AK1").Value = "pf_vm_c_to_pu"

    $ws.Range("AL1").Value = "pf_va_a_from_degree"
    $ws.Range("AM1").Value = "pf_va_b_from_degree"
    $ws.Range("AN1").Value = "pf_va_c_from_degree"
    $ws.Range("AO1").Value = "pf_va_a_to_degree"
    $ws.Range("AP1").Value = "pf_va_b_to_degree"
    $ws.Range("AQ1").Value = "pf_va_c_to_degree"
}
